# Update the division-fact table with a new set of problems.
# The document's single table has data rows at table-rows 1, 5, 9, 13, 17
# (1-based), each with 5 populated cells; the intervening rows are blank
# spacer rows. We address each cell by its (row, column) table position so
# that duplicate text values (e.g. "46÷9=") and values that collide with
# other rows' replacement text (e.g. "45÷6=") are never ambiguous, as a
# plain document-wide Find/Replace would be.

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Text = "42÷5=" },
    @{ Row = 1;  Col = 2; Text = "42÷9=" },
    @{ Row = 1;  Col = 3; Text = "45÷6=" },
    @{ Row = 1;  Col = 4; Text = "88÷4=" },
    @{ Row = 1;  Col = 5; Text = "69÷8=" },

    @{ Row = 5;  Col = 1; Text = "24÷7=" },
    @{ Row = 5;  Col = 2; Text = "67÷9=" },
    @{ Row = 5;  Col = 3; Text = "21÷8=" },
    @{ Row = 5;  Col = 4; Text = "52÷3=" },
    @{ Row = 5;  Col = 5; Text = "55÷7=" },

    @{ Row = 9;  Col = 1; Text = "45÷5=" },
    @{ Row = 9;  Col = 2; Text = "33÷7=" },
    @{ Row = 9;  Col = 3; Text = "85÷4=" },
    @{ Row = 9;  Col = 4; Text = "89÷6=" },
    @{ Row = 9;  Col = 5; Text = "59÷2=" },

    @{ Row = 13; Col = 1; Text = "85÷8=" },
    @{ Row = 13; Col = 2; Text = "38÷5=" },
    @{ Row = 13; Col = 3; Text = "63÷8=" },
    @{ Row = 13; Col = 4; Text = "41÷9=" },
    @{ Row = 13; Col = 5; Text = "44÷2=" },

    @{ Row = 17; Col = 1; Text = "41÷3=" },
    @{ Row = 17; Col = 2; Text = "21÷6=" },
    @{ Row = 17; Col = 3; Text = "69÷3=" },
    @{ Row = 17; Col = 4; Text = "68÷3=" },
    @{ Row = 17; Col = 5; Text = "82÷8=" }
)

foreach ($u in $updates) {
    $table.Cell($u.Row, $u.Col).Range.Text = $u.Text
}
